$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'305.53"
$ws.Range('E2').Value = "'0.94%"
$ws.Range('D3').Value = "'36.04"
$ws.Range('E3').Value = "'-3.12%"
$ws.Range('D4').Value = "'5.101"
$ws.Range('E4').Value = "'2.20%"
$ws.Range('D5').Value = "'0.07888"
$ws.Range('E5').Value = "'0.78%"
$ws.Range('D6').Value = "'2.143"
$ws.Range('E6').Value = "'-2.84%"
$ws.Range('D7').Value = "'7.929"
$ws.Range('E7').Value = "'-1.02%"
$ws.Range('D8').Value = "'0.9232"
$ws.Range('E8').Value = "'0.99%"
$ws.Range('D9').Value = "'0.09750"
$ws.Range('E9').Value = "'-0.08%"
$ws.Range('D10').Value = "'0.1859"
$ws.Range('E10').Value = "'-1.48%"
$ws.Range('D11').Value = "'0.08689"
$ws.Range('E11').Value = "'0.39%"
$ws.Range('D12').Value = "'0.03558"
$ws.Range('E12').Value = "'-0.31%"
$ws.Range('D13').Value = "'0.09943"
$ws.Range('E13').Value = "'-0.20%"
$ws.Range('D14').Value = "'0.001433"
$ws.Range('E14').Value = "'-3.57%"
$ws.Range('D15').Value = "'0.005628"
$ws.Range('E15').Value = "'-0.38%"
$ws.Range('D16').Value = "'3.455"
$ws.Range('E16').Value = "'-0.22%"
$ws.Range('D17').Value = "'4.107"
$ws.Range('E17').Value = "'1.93%"
$ws.Range('D18').Value = "'2.620"
$ws.Range('E18').Value = "'15.83%"
$ws.Range('D19').Value = "'0.3395"
$ws.Range('E19').Value = "'-1.93%"
$ws.Range('D20').Value = "'0.1318"
$ws.Range('E20').Value = "'1.36%"
$ws.Range('D21').Value = "'5.185"
$ws.Range('E21').Value = "'8.88%"
$ws.Range('D23').Value = "'0.04564"
$ws.Range('E23').Value = "'-1.53%"
$ws.Range('D24').Value = "'0.005053"
$ws.Range('E24').Value = "'5.55%"
$ws.Range('D25').Value = "'0.001234"
$ws.Range('E25').Value = "'0.15%"
$ws.Range('D27').Value = "'0.0004742"
$ws.Range('D39').Value = "'0.01855"
$ws.Range('E39').Value = "'4.58%"
$ws.Range('D40').Value = "'0.04764"
$ws.Range('E40').Value = "'0.50%"
$ws.Range('D41').Value = "'0.007501"
$ws.Range('E41').Value = "'-6.97%"
$ws.Range('D42').Value = "'0.1400"
$ws.Range('E42').Value = "'0.58%"
$ws.Range('D43').Value = "'0.007733"
$ws.Range('E43').Value = "'0.92%"
$ws.Range('D44').Value = "'0.002227"
$ws.Range('E44').Value = "'3.04%"
$ws.Range('D45').Value = "'0.01133"
$ws.Range('E45').Value = "'14.77%"
$ws.Range('D46').Value = "'0.00006321"
$ws.Range('E46').Value = "'3.57%"
$ws.Range('E47').Value = "'-0.18%"
$ws.Range('D48').Value = "'0.0005794"
$ws.Range('E48').Value = "'-0.11%"
$ws.Range('D49').Value = "'47.49"
$ws.Range('E49').Value = "'500.66%"
$ws.Range('E51').Value = "'-0.18%"
